$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "date" column (D) values, row 2..89, in order (replacing the shared-string
# content of the date column with updated timestamps).
$dates = @(
    "2018-07-31 14:53:24",
    "2018-07-31 14:53:36",
    "2018-07-31 14:53:48",
    "2018-07-31 14:53:55",
    "2018-07-31 14:54:07",
    "2018-07-31 14:54:18",
    "2018-07-31 14:54:27",
    "2018-07-31 14:54:38",
    "2018-07-31 14:54:48",
    "2018-07-31 14:54:59",
    "2018-07-31 14:55:12",
    "2018-07-31 14:55:19",
    "2018-07-31 14:55:30",
    "2018-07-31 14:55:40",
    "2018-07-31 14:55:50",
    "2018-07-31 14:56:01",
    "2018-07-31 14:56:06",
    "2018-07-31 14:56:19",
    "2018-07-31 14:56:29",
    "2018-07-31 14:56:36",
    "2018-07-31 14:56:46",
    "2018-07-31 14:56:56",
    "2018-07-31 14:57:06",
    "2018-07-31 14:57:16",
    "2018-07-31 14:57:23",
    "2018-07-31 14:57:34",
    "2018-07-31 14:57:42",
    "2018-07-31 14:57:49",
    "2018-07-31 14:58:01",
    "2018-07-31 14:58:11",
    "2018-07-31 14:58:22",
    "2018-07-31 14:58:33",
    "2018-07-31 14:58:40",
    "2018-07-31 14:58:51",
    "2018-07-31 14:59:03",
    "2018-07-31 14:59:14",
    "2018-07-31 14:59:23",
    "2018-07-31 14:59:31",
    "2018-07-31 14:59:42",
    "2018-07-31 14:59:52",
    "2018-07-31 15:00:01",
    "2018-07-31 15:00:06",
    "2018-07-31 15:00:18",
    "2018-07-31 15:00:30",
    "2018-07-31 15:00:37",
    "2018-07-31 15:00:48",
    "2018-07-31 15:00:59",
    "2018-07-31 15:01:12",
    "2018-07-31 15:01:23",
    "2018-07-31 15:01:34",
    "2018-07-31 15:01:44",
    "2018-07-31 15:01:54",
    "2018-07-31 15:02:03",
    "2018-07-31 15:02:12",
    "2018-07-31 15:02:22",
    "2018-07-31 15:02:32",
    "2018-07-31 15:02:41",
    "2018-07-31 15:02:51",
    "2018-07-31 15:03:01",
    "2018-07-31 15:03:10",
    "2018-07-31 15:03:22",
    "2018-07-31 15:03:33",
    "2018-07-31 15:03:44",
    "2018-07-31 15:03:56",
    "2018-07-31 15:04:02",
    "2018-07-31 15:04:14",
    "2018-07-31 15:04:25",
    "2018-07-31 15:04:33",
    "2018-07-31 15:04:45",
    "2018-07-31 15:04:57",
    "2018-07-31 15:05:05",
    "2018-07-31 15:05:20",
    "2018-07-31 15:05:32",
    "2018-07-31 15:05:42",
    "2018-07-31 15:05:54",
    "2018-07-31 15:06:06",
    "2018-07-31 15:06:16",
    "2018-07-31 15:06:24",
    "2018-07-31 15:06:34",
    "2018-07-31 15:06:45",
    "2018-07-31 15:06:57",
    "2018-07-31 15:07:11",
    "2018-07-31 15:07:20",
    "2018-07-31 15:07:30",
    "2018-07-31 15:07:40",
    "2018-07-31 15:07:47",
    "2018-07-31 15:07:57",
    "2018-07-31 15:08:09",
)

for ($i = 0; $i -lt $dates.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dates[$i]
}

# Updated det_int / det_int_redundant (columns B, C) values for the rows that
# changed between the two detector runs.
$bcChanges = @{
    2 = @(0, 0)
    3 = @(0, 0)
    4 = @(0, 0)
    5 = @(0, 0)
    6 = @(0, 0)
    17 = @(10, 10)
    18 = @(10, 10)
    19 = @(10, 10)
    20 = @(10, 10)
    21 = @(10, 10)
    22 = @(10, 10)
    23 = @(10, 10)
    24 = @(10, 10)
    25 = @(10, 10)
    26 = @(20, 20)
    27 = @(10, 10)
    28 = @(10, 10)
    29 = @(10, 10)
    30 = @(10, 10)
    31 = @(10, 10)
    33 = @(0, 0)
    39 = @(10, 10)
    40 = @(10, 10)
    41 = @(10, 10)
    42 = @(10, 10)
    43 = @(10, 10)
    50 = @(10, 10)
    51 = @(10, 10)
    52 = @(10, 10)
    53 = @(10, 10)
    54 = @(10, 10)
    55 = @(10, 10)
    78 = @(10, 10)
    81 = @(0, 0)
    82 = @(0, 0)
    83 = @(0, 0)
    85 = @(0, 0)
    87 = @(0, 0)
}

foreach ($row in $bcChanges.Keys) {
    $pair = $bcChanges[$row]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
}

# Reset the view: scroll back to the top and select A1 instead of the
# previously scrolled/selected cell (D24, topLeftCell A55).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
